$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New helper column J: flag (1/0) rows whose effort exceeds 3 ---
# Functional estimate table (rows 8-21)
$ws.Range("J8").Formula = "=IF(E8>3,1,0)"

# Technical estimate table (rows 28-36)
$ws.Range("J28").Formula = "=IF(F28>3,1,0)"
$ws.Range("J29:J35").Formula = "=IF(F29>3,1,0)"
$ws.Range("J36").Formula = "=SUM(J28:J35)"

# --- New columns L (effort offset) and M (running total) for rows 29-31 ---
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 20

$ws.Range("L30").Value = 3
$ws.Range("M30").Value = 28

$ws.Range("L31").Value = 5
$ws.Range("M31").Value = 36

# --- View / selection updates ---
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("L36").Select()
